$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.587.72'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.648.58'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.69'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('E6').Value = '  +4.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.64'
$ws.Range('E8').Value = '  -2.92%  '
$ws.Range('E9').Value = '  -1.63%  '
$ws.Range('E10').Value = '  -1.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0892'
$ws.Range('D12').Value = '1.880.63'
$ws.Range('E12').Value = '  -0.96%  '
$ws.Range('D13').Value = '1.659.41'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.587'
$ws.Range('E14').Value = '  +3.42%  '
$ws.Range('E15').Value = '  -2.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.52'
$ws.Range('E16').Value = '  -2.45%  '
$ws.Range('D17').Value = '27.557.42'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.48'
$ws.Range('E18').Value = '  -3.68%  '
$ws.Range('D19').Value = '0.0₃0727'
$ws.Range('E20').Value = '  -0.92%  '
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.34'
$ws.Range('E22').Value = '  -3.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.76'
$ws.Range('E23').Value = '  +3.91%  '
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.49'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.04'
$ws.Range('E26').Value = '  -2.81%  '
$ws.Range('E27').Value = '  +1.71%  '
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.63'
$ws.Range('E29').Value = '  -4.42%  '
$ws.Range('E30').Value = '  -2.38%  '
$ws.Range('E31').Value = '  -3.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.32'
$ws.Range('E32').Value = '  -0.92%  '
$ws.Range('E33').Value = '  +1.45%  '
$ws.Range('D34').Value = '1.426.31'
$ws.Range('E34').Value = '  -2.46%  '
$ws.Range('E35').Value = '  +1.61%  '
$ws.Range('E36').Value = '  -0.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.570'
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('E38').Value = '  -4.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0166'
$ws.Range('E39').Value = '  -3.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.02'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.819'
$ws.Range('E42').Value = '  +3.42%  '
$ws.Range('E43').Value = '  -2.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.55'
$ws.Range('E44').Value = '  +2.50%  '
$ws.Range('E45').Value = '  +1.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.20'
$ws.Range('E46').Value = '  -7.10%  '
$ws.Range('D47').Value = '1.791.22'
$ws.Range('E47').Value = '  -0.84%  '
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '88.15'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').Value = '0.0₆0107'
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.80'
$ws.Range('E51').Value = '  -0.99%  '
